# Generate Report for Handback
# Updates the two "handed back" file records (uuid-named .md sources and
# their related .xlf target files) with freshly generated identifiers /
# timestamps, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$wsOverview.Range("B2").Value = "e2e\6a1f3617-07d8-4b4c-add0-535573617591.md"
$wsOverview.Range("G2").Value = "2016-08-29 07:04:05"

$wsOverview.Range("A3").Value = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$wsOverview.Range("B3").Value = "e2e\ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$wsOverview.Range("G3").Value = "2016-08-29 07:04:05"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\6a1f3617-07d8-4b4c-add0-535573617591.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$wsZhCn.Range("G2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-29 07:03:57"
$wsZhCn.Range("I2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$wsZhCn.Range("J2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-29 07:04:30"

$wsZhCn.Range("A3").Value = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$wsZhCn.Range("G3").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-29 07:03:57"
$wsZhCn.Range("I3").Value = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$wsZhCn.Range("J3").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-29 07:04:30"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "6a1f3617-07d8-4b4c-add0-535573617591.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$wsDeDe.Range("G2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-29 07:04:05"
$wsDeDe.Range("I2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$wsDeDe.Range("J2").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-29 07:04:37"

$wsDeDe.Range("A3").Value = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$wsDeDe.Range("G3").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-29 07:04:05"
$wsDeDe.Range("I3").Value = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
$wsDeDe.Range("J3").Value = "6a1f3617-07d8-4b4c-add0-535573617591.21614cdb91720701de0327b3fbcce43157ad5500.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-29 07:04:37"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "6a1f3617-07d8-4b4c-add0-535573617591.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "ffff09a99ca6-b8cb-4b53-83b7-9f1935683d43.md"
    }
}

Write-Output "Handback report regenerated"
